$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "heal-ball",
    "burn-heal",
    "ice-heal",
    "paralyze-heal",
    "full-heal",
    "heal-powder",
    "health-wing",
    "health-candy",
    "health-candy-l"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

# Remove the now-unused rows (previously rows 11-35)
$ws.Range("A11:A35").ClearContents() | Out-Null
